$wb = $excel.ActiveWorkbook

# Sheet "展览" (first sheet) - column F updates
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 6832
$ws1.Range("F4").Value = 22
$ws1.Range("F6").Value = 0
$ws1.Range("F7").Value = 6498
$ws1.Range("F8").Value = 0
$ws1.Range("F10").Value = 1286
$ws1.Range("F11").Value = 0
$ws1.Range("F13").Value = 0
$ws1.Range("F14").Value = 0
$ws1.Range("F15").Value = 0
$ws1.Range("F16").Value = 380
$ws1.Range("F17").Value = 0
$ws1.Range("F18").Value = 0
$ws1.Range("F19").Value = 4841
$ws1.Range("F21").Value = 77
$ws1.Range("F22").Value = 0
$ws1.Range("F23").Value = 0

# Sheet "全部类型" (fourth sheet) - column F updates
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 0
$ws4.Range("F4").Value = 0
$ws4.Range("F7").Value = 0
$ws4.Range("F8").Value = 0
$ws4.Range("F9").Value = 0
$ws4.Range("F10").Value = 0
$ws4.Range("F11").Value = 0
$ws4.Range("F13").Value = 395
$ws4.Range("F15").Value = 0
$ws4.Range("F16").Value = 380
$ws4.Range("F18").Value = 0
$ws4.Range("F20").Value = 0
$ws4.Range("F22").Value = 0
$ws4.Range("F23").Value = 305
$ws4.Range("F25").Value = 0

$wb.Save()
